# The deck's two theme parts (ppt/theme/theme1.xml, used by the slide
# master, and ppt/theme/theme2.xml, used by the notes master) swap
# content: theme1 becomes the stock "Office Theme" palette (previously
# theme2's content) and theme2 becomes the old "Integral" palette
# (previously theme1's content).
#
# The font scheme (Arial-based "Office" fonts) and the format scheme
# (fill/line/effect styles) are already byte-identical between the two
# theme parts, so only the 12 color-scheme slots (plus the cosmetic
# "name" attributes, which PowerPoint does not expose as writable
# anywhere in the object model - ThemeColorScheme.Name is documented
# read-only) actually need to change.
#
# Apply the new ("Office Theme") colors to the presentation's theme
# color scheme, in MsoThemeColorSchemeIndex order:
#   1 Dark1, 2 Light1, 3 Dark2, 4 Light2,
#   5-10 Accent1-6, 11 Hyperlink, 12 FollowedHyperlink
# (PowerPoint's RGB value is 0x00BBGGRR, i.e. blue/green/red packed low-to-high.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # Dark 1    - 000000
$tcs.Item(2).RGB  = 16777215   # Light 1   - FFFFFF
$tcs.Item(3).RGB  = 6968388    # Dark 2    - 44546A
$tcs.Item(4).RGB  = 15132391   # Light 2   - E7E6E6
$tcs.Item(5).RGB  = 13998939   # Accent 1  - 5B9BD5
$tcs.Item(6).RGB  = 3243501    # Accent 2  - ED7D31
$tcs.Item(7).RGB  = 10855845   # Accent 3  - A5A5A5
$tcs.Item(8).RGB  = 49407      # Accent 4  - FFC000
$tcs.Item(9).RGB  = 12874308   # Accent 5  - 4472C4
$tcs.Item(10).RGB = 4697456    # Accent 6  - 70AD47
$tcs.Item(11).RGB = 12673797   # Hyperlink - 0563C1
$tcs.Item(12).RGB = 7491477    # Followed Hyperlink - 954F72
